$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G, shifting existing G:J to H:K
$ws.Range("G1").EntireColumn.Insert()

# Set the new header for column G
$ws.Range("G1").Value = "sum_assured"

# Fill the new sum_assured column with a constant value for all data rows
$ws.Range("G2:G5").Value = 500000
